# Updated code Date 6 and 6PM
# - Update the email address in C2 (jerry3 -> jerry4 @yopmail.com)
# - Turn that cell into a mailto hyperlink (adds the built-in "Hyperlink" style)
# - Move the active selection to F2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "jerry4@yopmail.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:jerry4@yopmail.com") | Out-Null
$ws.Range("F2").Select() | Out-Null
